$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9
$ws.Cells.Item($row, 1).Value = 42612.892488425925
$ws.Cells.Item($row, 2).Value = 20
$ws.Cells.Item($row, 3).Value = 54
$ws.Cells.Item($row, 4).Value = 41
$ws.Cells.Item($row, 5).Value = 82
$ws.Cells.Item($row, 6).Value = 17
$ws.Cells.Item($row, 7).Value = 13103
$ws.Cells.Item($row, 8).Value = 26404
$ws.Cells.Item($row, 9).Value = 3223
$ws.Cells.Item($row, 10).Value = 364
$ws.Cells.Item($row, 11).Value = 279
$ws.Cells.Item($row, 12).Value = 24
$ws.Cells.Item($row, 13).Value = 5
$ws.Cells.Item($row, 14).Value = "Bag"

$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"
